$wb = $excel.ActiveWorkbook

# --- Sheet 1 (БИВТ-22-17) ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 5

$ws1.Range("C3").Value = "pass"

$ws1.Range("D5").Value = 5

$ws1.Range("C9").Value = 5
$ws1.Range("D9").Value = 5

$ws1.Range("C11").Value = 5

$ws1.Range("D17").Value = 5

$ws1.Range("D19").Value = 5

$ws1.Range("C20").Value = 5

$ws1.Range("D23").Value = 5
$ws1.Range("E23").Value = 5

$ws1.Range("C26").Value = 5

$ws1.Range("D28").Value = 5

# --- Sheet 2 (БИВТ-22-18) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C2").Value = "pass"
